$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date style (s="1", numFmt 14) from B743 down to the new date cells so we reuse the existing style index
$ws.Cells.Item(743, 2).Copy($ws.Range("B744:B757"))

# Step 1: fill column G (Temps joue) first, top-to-bottom, to match the shared-string insertion order in the target diff
$ws.Cells.Item(744, 7).Value = "01:06:30"
$ws.Cells.Item(745, 7).Value = "01:06:06"
$ws.Cells.Item(746, 7).Value = "00:26:14"
$ws.Cells.Item(747, 7).Value = "00:25:58"
$ws.Cells.Item(748, 7).Value = "01:31:56"
$ws.Cells.Item(749, 7).Value = "01:32:28"
$ws.Cells.Item(750, 7).Value = "01:32:28"
$ws.Cells.Item(751, 7).Value = "00:26:22"
$ws.Cells.Item(752, 7).Value = "01:32:44"
$ws.Cells.Item(753, 7).Value = "00:46:46"
$ws.Cells.Item(754, 7).Value = "01:32:28"
$ws.Cells.Item(755, 7).Value = "01:32:28"
$ws.Cells.Item(756, 7).Value = "00:45:42"
$ws.Cells.Item(757, 7).Value = "01:05:42"

# Step 2: fill the rest of each row (B, C, E, F, H..V), leaving column A (Type) for last
$ws.Cells.Item(744, 2).Value = 45955
$ws.Cells.Item(744, 3).Value = "Global"
$ws.Cells.Item(744, 5).Value = "Naim Dhib"
$ws.Cells.Item(744, 6).Value = "center midfield"
$ws.Cells.Item(744, 8).Value = 7.72
$ws.Cells.Item(744, 9).Value = 1.26
$ws.Cells.Item(744, 10).Value = 6.45
$ws.Cells.Item(744, 11).Value = 0.85
$ws.Cells.Item(744, 12).Value = 0.36
$ws.Cells.Item(744, 13).Value = 0.06
$ws.Cells.Item(744, 14).Value = 0
$ws.Cells.Item(744, 15).Value = 7
$ws.Cells.Item(744, 16).Value = 6.93
$ws.Cells.Item(744, 17).Value = 27.79
$ws.Cells.Item(744, 18).Value = 4.5
$ws.Cells.Item(744, 19).Value = 32
$ws.Cells.Item(744, 20).Value = 2
$ws.Cells.Item(744, 21).Value = 32
$ws.Cells.Item(744, 22).Value = 7
$ws.Cells.Item(745, 2).Value = 45955
$ws.Cells.Item(745, 3).Value = "Global"
$ws.Cells.Item(745, 5).Value = "Malik Boussaid"
$ws.Cells.Item(745, 6).Value = "right back"
$ws.Cells.Item(745, 8).Value = 8.5399999999999991
$ws.Cells.Item(745, 9).Value = 2.08
$ws.Cells.Item(745, 10).Value = 6.44
$ws.Cells.Item(745, 11).Value = 1.22
$ws.Cells.Item(745, 12).Value = 0.65
$ws.Cells.Item(745, 13).Value = 0.22
$ws.Cells.Item(745, 14).Value = 0.02
$ws.Cells.Item(745, 15).Value = 12
$ws.Cells.Item(745, 16).Value = 7.72
$ws.Cells.Item(745, 17).Value = 31.89
$ws.Cells.Item(745, 18).Value = 4.46
$ws.Cells.Item(745, 19).Value = 49
$ws.Cells.Item(745, 20).Value = 5
$ws.Cells.Item(745, 21).Value = 23
$ws.Cells.Item(745, 22).Value = 10
$ws.Cells.Item(746, 2).Value = 45955
$ws.Cells.Item(746, 3).Value = "Global"
$ws.Cells.Item(746, 5).Value = "Karahali Souaré"
$ws.Cells.Item(746, 6).Value = "right forward"
$ws.Cells.Item(746, 8).Value = 3.2
$ws.Cells.Item(746, 9).Value = 0.73
$ws.Cells.Item(746, 10).Value = 2.46
$ws.Cells.Item(746, 11).Value = 0.4
$ws.Cells.Item(746, 12).Value = 0.22
$ws.Cells.Item(746, 13).Value = 0.12
$ws.Cells.Item(746, 14).Value = 0
$ws.Cells.Item(746, 15).Value = 8
$ws.Cells.Item(746, 16).Value = 7.3
$ws.Cells.Item(746, 17).Value = 29.85
$ws.Cells.Item(746, 18).Value = 4.6399999999999997
$ws.Cells.Item(746, 19).Value = 13
$ws.Cells.Item(746, 20).Value = 7
$ws.Cells.Item(746, 21).Value = 10
$ws.Cells.Item(746, 22).Value = 9
$ws.Cells.Item(747, 2).Value = 45955
$ws.Cells.Item(747, 3).Value = "Global"
$ws.Cells.Item(747, 5).Value = "Amir Etien"
$ws.Cells.Item(747, 6).Value = "right forward"
$ws.Cells.Item(747, 8).Value = 2.61
$ws.Cells.Item(747, 9).Value = 0.57999999999999996
$ws.Cells.Item(747, 10).Value = 2.0299999999999998
$ws.Cells.Item(747, 11).Value = 0.35
$ws.Cells.Item(747, 12).Value = 0.2
$ws.Cells.Item(747, 13).Value = 0.02
$ws.Cells.Item(747, 14).Value = 0.01
$ws.Cells.Item(747, 15).Value = 2
$ws.Cells.Item(747, 16).Value = 6.04
$ws.Cells.Item(747, 17).Value = 31.31
$ws.Cells.Item(747, 18).Value = 4.99
$ws.Cells.Item(747, 19).Value = 12
$ws.Cells.Item(747, 20).Value = 5
$ws.Cells.Item(747, 21).Value = 4
$ws.Cells.Item(747, 22).Value = 2
$ws.Cells.Item(748, 2).Value = 45955
$ws.Cells.Item(748, 3).Value = "Global"
$ws.Cells.Item(748, 5).Value = "Yoann Martelat"
$ws.Cells.Item(748, 6).Value = "center midfield"
$ws.Cells.Item(748, 8).Value = 11.92
$ws.Cells.Item(748, 9).Value = 2.16
$ws.Cells.Item(748, 10).Value = 9.73
$ws.Cells.Item(748, 11).Value = 1.76
$ws.Cells.Item(748, 12).Value = 0.39
$ws.Cells.Item(748, 13).Value = 0.04
$ws.Cells.Item(748, 14).Value = 0
$ws.Cells.Item(748, 15).Value = 6
$ws.Cells.Item(748, 16).Value = 7.76
$ws.Cells.Item(748, 17).Value = 26.92
$ws.Cells.Item(748, 18).Value = 4.13
$ws.Cells.Item(748, 19).Value = 24
$ws.Cells.Item(748, 20).Value = 1
$ws.Cells.Item(748, 21).Value = 27
$ws.Cells.Item(748, 22).Value = 1
$ws.Cells.Item(749, 2).Value = 45955
$ws.Cells.Item(749, 3).Value = "Global"
$ws.Cells.Item(749, 5).Value = "Naim Ighbane"
$ws.Cells.Item(749, 6).Value = "center back"
$ws.Cells.Item(749, 8).Value = 9.48
$ws.Cells.Item(749, 9).Value = 1.27
$ws.Cells.Item(749, 10).Value = 8.1999999999999993
$ws.Cells.Item(749, 11).Value = 0.91
$ws.Cells.Item(749, 12).Value = 0.28000000000000003
$ws.Cells.Item(749, 13).Value = 0.08
$ws.Cells.Item(749, 14).Value = 0
$ws.Cells.Item(749, 15).Value = 8
$ws.Cells.Item(749, 16).Value = 6.14
$ws.Cells.Item(749, 17).Value = 29.71
$ws.Cells.Item(749, 18).Value = 4.49
$ws.Cells.Item(749, 19).Value = 32
$ws.Cells.Item(749, 20).Value = 3
$ws.Cells.Item(749, 21).Value = 22
$ws.Cells.Item(749, 22).Value = 3
$ws.Cells.Item(750, 2).Value = 45955
$ws.Cells.Item(750, 3).Value = "Global"
$ws.Cells.Item(750, 5).Value = "Sofiane Belle"
$ws.Cells.Item(750, 6).Value = "left forward"
$ws.Cells.Item(750, 8).Value = 10.220000000000001
$ws.Cells.Item(750, 9).Value = 1.91
$ws.Cells.Item(750, 10).Value = 8.2799999999999994
$ws.Cells.Item(750, 11).Value = 1.23
$ws.Cells.Item(750, 12).Value = 0.52
$ws.Cells.Item(750, 13).Value = 0.19
$ws.Cells.Item(750, 14).Value = 0
$ws.Cells.Item(750, 15).Value = 11
$ws.Cells.Item(750, 16).Value = 6.64
$ws.Cells.Item(750, 17).Value = 29.7
$ws.Cells.Item(750, 18).Value = 4.1900000000000004
$ws.Cells.Item(750, 19).Value = 27
$ws.Cells.Item(750, 20).Value = 1
$ws.Cells.Item(750, 21).Value = 26
$ws.Cells.Item(750, 22).Value = 9
$ws.Cells.Item(751, 2).Value = 45955
$ws.Cells.Item(751, 3).Value = "Global"
$ws.Cells.Item(751, 5).Value = "Mattheo Haon"
$ws.Cells.Item(751, 6).Value = "right back"
$ws.Cells.Item(751, 8).Value = 3.15
$ws.Cells.Item(751, 9).Value = 0.85
$ws.Cells.Item(751, 10).Value = 2.2999999999999998
$ws.Cells.Item(751, 11).Value = 0.51
$ws.Cells.Item(751, 12).Value = 0.26
$ws.Cells.Item(751, 13).Value = 0.08
$ws.Cells.Item(751, 14).Value = 0
$ws.Cells.Item(751, 15).Value = 7
$ws.Cells.Item(751, 16).Value = 7.14
$ws.Cells.Item(751, 17).Value = 28.25
$ws.Cells.Item(751, 18).Value = 5.61
$ws.Cells.Item(751, 19).Value = 6
$ws.Cells.Item(751, 20).Value = 5
$ws.Cells.Item(751, 21).Value = 9
$ws.Cells.Item(751, 22).Value = 6
$ws.Cells.Item(752, 2).Value = 45955
$ws.Cells.Item(752, 3).Value = "Global"
$ws.Cells.Item(752, 5).Value = "Romain Thunet"
$ws.Cells.Item(752, 6).Value = "center back"
$ws.Cells.Item(752, 8).Value = 10.82
$ws.Cells.Item(752, 9).Value = 1.79
$ws.Cells.Item(752, 10).Value = 9
$ws.Cells.Item(752, 11).Value = 1.28
$ws.Cells.Item(752, 12).Value = 0.43
$ws.Cells.Item(752, 13).Value = 0.1
$ws.Cells.Item(752, 14).Value = 0
$ws.Cells.Item(752, 15).Value = 7
$ws.Cells.Item(752, 16).Value = 7.01
$ws.Cells.Item(752, 17).Value = 27.73
$ws.Cells.Item(752, 18).Value = 4.62
$ws.Cells.Item(752, 19).Value = 50
$ws.Cells.Item(752, 20).Value = 5
$ws.Cells.Item(752, 21).Value = 41
$ws.Cells.Item(752, 22).Value = 4
$ws.Cells.Item(753, 2).Value = 45955
$ws.Cells.Item(753, 3).Value = "Global"
$ws.Cells.Item(753, 5).Value = "Emmanuel Valey"
$ws.Cells.Item(753, 6).Value = "left forward"
$ws.Cells.Item(753, 8).Value = 5.86
$ws.Cells.Item(753, 9).Value = 1.45
$ws.Cells.Item(753, 10).Value = 4.4000000000000004
$ws.Cells.Item(753, 11).Value = 0.8
$ws.Cells.Item(753, 12).Value = 0.41
$ws.Cells.Item(753, 13).Value = 0.22
$ws.Cells.Item(753, 14).Value = 0.03
$ws.Cells.Item(753, 15).Value = 18
$ws.Cells.Item(753, 16).Value = 7.51
$ws.Cells.Item(753, 17).Value = 31.7
$ws.Cells.Item(753, 18).Value = 5.22
$ws.Cells.Item(753, 19).Value = 34
$ws.Cells.Item(753, 20).Value = 9
$ws.Cells.Item(753, 21).Value = 25
$ws.Cells.Item(753, 22).Value = 13
$ws.Cells.Item(754, 2).Value = 45955
$ws.Cells.Item(754, 3).Value = "Global"
$ws.Cells.Item(754, 5).Value = "Ilyes Boughanmi"
$ws.Cells.Item(754, 6).Value = "center forward"
$ws.Cells.Item(754, 8).Value = 9.86
$ws.Cells.Item(754, 9).Value = 2.13
$ws.Cells.Item(754, 10).Value = 7.71
$ws.Cells.Item(754, 11).Value = 1.32
$ws.Cells.Item(754, 12).Value = 0.61
$ws.Cells.Item(754, 13).Value = 0.2
$ws.Cells.Item(754, 14).Value = 0.03
$ws.Cells.Item(754, 15).Value = 18
$ws.Cells.Item(754, 16).Value = 6.38
$ws.Cells.Item(754, 17).Value = 31.12
$ws.Cells.Item(754, 18).Value = 4.49
$ws.Cells.Item(754, 19).Value = 46
$ws.Cells.Item(754, 20).Value = 5
$ws.Cells.Item(754, 21).Value = 55
$ws.Cells.Item(754, 22).Value = 13
$ws.Cells.Item(755, 2).Value = 45955
$ws.Cells.Item(755, 3).Value = "Global"
$ws.Cells.Item(755, 5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(755, 6).Value = "center midfield"
$ws.Cells.Item(755, 8).Value = 12.37
$ws.Cells.Item(755, 9).Value = 2.62
$ws.Cells.Item(755, 10).Value = 9.73
$ws.Cells.Item(755, 11).Value = 2.09
$ws.Cells.Item(755, 12).Value = 0.5
$ws.Cells.Item(755, 13).Value = 0.06
$ws.Cells.Item(755, 14).Value = 0
$ws.Cells.Item(755, 15).Value = 7
$ws.Cells.Item(755, 16).Value = 8.0299999999999994
$ws.Cells.Item(755, 17).Value = 26.89
$ws.Cells.Item(755, 18).Value = 4.74
$ws.Cells.Item(755, 19).Value = 49
$ws.Cells.Item(755, 20).Value = 10
$ws.Cells.Item(755, 21).Value = 26
$ws.Cells.Item(755, 22).Value = 4
$ws.Cells.Item(756, 2).Value = 45955
$ws.Cells.Item(756, 3).Value = "Global"
$ws.Cells.Item(756, 5).Value = "Maé Clavel"
$ws.Cells.Item(756, 6).Value = "left back"
$ws.Cells.Item(756, 8).Value = 5.25
$ws.Cells.Item(756, 9).Value = 1.1000000000000001
$ws.Cells.Item(756, 10).Value = 4.1399999999999997
$ws.Cells.Item(756, 11).Value = 0.72
$ws.Cells.Item(756, 12).Value = 0.25
$ws.Cells.Item(756, 13).Value = 0.15
$ws.Cells.Item(756, 14).Value = 0
$ws.Cells.Item(756, 15).Value = 8
$ws.Cells.Item(756, 16).Value = 6.94
$ws.Cells.Item(756, 17).Value = 29.44
$ws.Cells.Item(756, 18).Value = 4.5
$ws.Cells.Item(756, 19).Value = 25
$ws.Cells.Item(756, 20).Value = 2
$ws.Cells.Item(756, 21).Value = 22
$ws.Cells.Item(756, 22).Value = 7
$ws.Cells.Item(757, 2).Value = 45955
$ws.Cells.Item(757, 3).Value = "Global"
$ws.Cells.Item(757, 5).Value = "Levy Ndoutoume"
$ws.Cells.Item(757, 6).Value = "left back"
$ws.Cells.Item(757, 8).Value = 7.73
$ws.Cells.Item(757, 9).Value = 1.56
$ws.Cells.Item(757, 10).Value = 6.15
$ws.Cells.Item(757, 11).Value = 0.86
$ws.Cells.Item(757, 12).Value = 0.46
$ws.Cells.Item(757, 13).Value = 0.22
$ws.Cells.Item(757, 14).Value = 0.04
$ws.Cells.Item(757, 15).Value = 12
$ws.Cells.Item(757, 16).Value = 7.02
$ws.Cells.Item(757, 17).Value = 29.71
$ws.Cells.Item(757, 18).Value = 4.49
$ws.Cells.Item(757, 19).Value = 37
$ws.Cells.Item(757, 20).Value = 7
$ws.Cells.Item(757, 21).Value = 26
$ws.Cells.Item(757, 22).Value = 15

# Step 3: fill column A (Type) last, so its new shared string is appended after all the time strings
$ws.Cells.Item(744, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(745, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(746, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(747, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(748, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(749, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(750, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(751, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(752, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(753, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(754, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(755, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(756, 1).Value = "CDF T6 VS Revermontoise (R3)"
$ws.Cells.Item(757, 1).Value = "CDF T6 VS Revermontoise (R3)"

# Update the active selection to match the post-edit state
$ws.Range("D746").Select()
